# ---------------------------------------------------------------------------
# Edit summary (from the canonical OOXML diff):
#
# 1. Slide 5 contains a table (the 2nd shape) whose <a:tblPr><a:tableStyleId>
#    changes from {F56E307F-92A4-435E-A67D-44FC74AF9427} to
#    {CD2607A9-D8AA-4EBC-90DE-70E2ED2EF2D8} - i.e. a different built-in
#    PowerPoint table style was applied to the table via the Table Design
#    gallery.
#
# 2. The presentation's main theme (ppt/theme/theme1.xml, used by the
#    slide master) had its colour scheme changed from the custom
#    "Red Violet" palette (part of the "Integral" theme) to the stock
#    Office colour palette. Font scheme / format scheme (fills, lines,
#    effects) are identical between the two themes, only the 12 theme
#    colours actually change - this is what happens when the "Office"
#    theme colours are applied from the Design > Variants > Colors
#    gallery.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{CD2607A9-D8AA-4EBC-90DE-70E2ED2EF2D8}")

# --- 2. Switch the presentation's theme colours to the standard "Office" --
#        palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink in order).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$anySlide = $p.Slides.Item(1)
$themeColors = $anySlide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
